$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.277.39'
$ws.Range("E2").Value = '  +0.66%  '

$ws.Range("D3").Value = '1.666.46'
$ws.Range("E3").Value = '  +0.69%  '

$ws.Range("E4").Value = '  +0.74%  '

$ws.Range("E5").Value = '  +0.41%  '

$ws.Range("E6").Value = '  +1.81%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2639'
$ws.Range("E8").Value = '  +1.24%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06383'
$ws.Range("E9").Value = '  +0.55%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.55'
$ws.Range("E10").Value = '  +0.89%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07819'
$ws.Range("E11").Value = '  +0.39%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.561'
$ws.Range("E12").Value = '  +1.30%  '

$ws.Range("D13").Value = '1.670.31'
$ws.Range("E13").Value = '  +0.56%  '

$ws.Range("D14").Value = '1.895.56'
$ws.Range("E14").Value = '  +0.65%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5530'
$ws.Range("E15").Value = '  +1.03%  '

$ws.Range("D16").Value = '0.0₅8198'
$ws.Range("E16").Value = '  -0.11%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '65.73'
$ws.Range("E17").Value = '  +0.55%  '

$ws.Range("E18").Value = '  +0.77%  '

$ws.Range("E19").Value = '  +2.18%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '193.49'
$ws.Range("E20").Value = '  +1.10%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.22'
$ws.Range("E21").Value = '  +1.54%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.033'
$ws.Range("E22").Value = '  -0.02%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.011'
$ws.Range("E23").Value = '  +0.73%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '145.68'
$ws.Range("E24").Value = '  +2.48%  '

$ws.Range("E25").Value = '  -1.20%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.197'
$ws.Range("E26").Value = '  -0.70%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.12'
$ws.Range("E27").Value = '  -0.01%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.483'
$ws.Range("E28").Value = '  +3.73%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.05894'
$ws.Range("E29").Value = '  -0.12%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.280'
$ws.Range("E30").Value = '  -0.17%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.599'
$ws.Range("E31").Value = '  +2.08%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.276'
$ws.Range("E32").Value = '  +0.83%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.614'
$ws.Range("E33").Value = '  +1.40%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.9636'
$ws.Range("E34").Value = '  +0.99%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.826'
$ws.Range("E35").Value = '  +1.33%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.424'
$ws.Range("E36").Value = '  +0.55%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.5799'
$ws.Range("E37").Value = '  +1.81%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01606'
$ws.Range("E38").Value = '  -0.73%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.8642'
$ws.Range("E39").Value = '  +1.68%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.839'
$ws.Range("E40").Value = '  +0.31%  '

$ws.Range("B41").Value = 'PaxDollar'
$ws.Range("C41").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.009'
$ws.Range("E41").Value = '  +0.65%  '

$ws.Range("B42").Value = 'Maker'
$ws.Range("C42").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D42").Value = '1.047.95'
$ws.Range("E42").Value = '  +1.39%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '103.97'
$ws.Range("E43").Value = '  +0.81%  '

$ws.Range("D44").Value = '1.805.77'
$ws.Range("E44").Value = '  +0.43%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '57.71'
$ws.Range("E45").Value = '  +0.86%  '

$ws.Range("B46").Value = 'Frax'
$ws.Range("C46").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.013'
$ws.Range("E46").Value = '  +0.56%  '

$ws.Range("B47").Value = 'BabyDogeCoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D47").Value = '0.0₈105'
$ws.Range("E47").Value = '  -5.61%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4380'
$ws.Range("E48").Value = '  +1.85%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.034'
$ws.Range("E49").Value = '  +2.42%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05161'
$ws.Range("E50").Value = '  -0.15%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.426'
$ws.Range("E51").Value = '  -3.31%  '
